# Sync helper workbook update:
#  - Convert the B3:C54 metadata range into a proper Excel Table (ListObject)
#  - Sort the table data by the "Metadata Name (code)" column (ascending)
#    so comparisons against the generated JSON are stable
#  - Apply the TableStyleMedium6 table style
#  - Leave the selection on the first data cell (B10) instead of the
#    previously-selected bottom-of-range cell

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The metadata name/label pairs currently live in a plain range - turn it
# into a real table so it can be referenced / filtered / sorted reliably.
$dataRange = $ws.Range("B3:C54")
$tbl = $ws.ListObjects.Add(
    [Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange,
    $dataRange,
    $null,
    [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes
)
$tbl.Name = "Table1"
$tbl.TableStyle = "TableStyleMedium6"

# Sort the table rows by the metadata code column (column B) ascending,
# keeping the header row in place.
$tbl.Sort.SortFields.Clear()
$tbl.Sort.SortFields.Add(
    $ws.Range("B4:B54"),
    [Microsoft.Office.Interop.Excel.XlSortOn]::xlSortOnValues,
    [Microsoft.Office.Interop.Excel.XlSortOrder]::xlAscending
) | Out-Null
$tbl.Sort.Header = [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes
$tbl.Sort.Orientation = [Microsoft.Office.Interop.Excel.XlSortOrientation]::xlSortColumns
$tbl.Sort.Apply()

# Move the selection/scroll position up to the top of the (now sorted) data
# instead of leaving it parked on the old last-row cell.
$ws.Range("B10").Select()
